# Weekly "Forward Look" stats update - 08.08.25
#
# 1) Refresh the "as at" date in the intro paragraph (A2): 01 August 2025 -> 08 August 2025
# 2) Remove the "Tribunals statistics quarterly: April to June 2025" publication row
#    (week commencing 08 Sep 2025) - clears Publication Title/Date/Status/Type but keeps
#    the week-commencing date and week number in place (same as other "no publication"
#    weeks in the sheet).
# 3) Widen the scope of the "Tribunals statistics quarterly" publication due 11 Dec 2025
#    from "July to September 2025" to "April to September 2025".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forward Look")

$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 08 August 2025"

$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("B30").Value = "Tribunals statistics quarterly: April to September 2025"
